# Add data for 2022-06-12: update the "through" date from 06-03 to 06-04
# and add the new day's carjacking counts for June across all years.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet and update its tab/title to reflect the new "through" date.
$ws.Name = "Through 2022-06-04"

# Update the "June (through 06-03)" label cell to "June (through 06-04)".
$ws.Range("A7").Value = "June (through 06-04)"

# Update June row (row 7) values for each year column (B..I).
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 10
$ws.Range("E7").Value = 12
$ws.Range("F7").Value = 5
$ws.Range("G7").Value = 29
$ws.Range("H7").Value = 15
$ws.Range("I7").Value = 10

# Update Total row (row 8) values for each year column (B..I).
$ws.Range("C8").Value = 214
$ws.Range("D8").Value = 326
$ws.Range("E8").Value = 307
$ws.Range("F8").Value = 209
$ws.Range("G8").Value = 387
$ws.Range("H8").Value = 646
$ws.Range("I8").Value = 674
